$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before the existing row 792 ("2026/12/29" / 火 / 13),
# shifting it (and every row below) down by one. The new row becomes
# "2026/02/09" / 月 / 19 / 201, matching the commit's daily push date.
$ws.Rows.Item(792).Insert()

# Column A holds the date stored as plain text (not an Excel date serial) in
# this workbook, so force a text format before assigning the value, then
# drop the format again so the cell ends up with no style, same as its
# neighbours.
$ws.Cells.Item(792, 1).NumberFormat = "@"
$ws.Cells.Item(792, 1).Value = "2026/02/09"
$ws.Cells.Item(792, 1).ClearFormats()

$ws.Cells.Item(792, 2).Value = "月"
$ws.Cells.Item(792, 3).Value = 19
$ws.Cells.Item(792, 4).Value = 201
